# System Setup4: sydefault, Inv comcod
# Edits the "Edit_Order Number Assign" sheet: bump the default/starting
# sequence numbers, re-fit the column widths to the (now longer) header
# text, and leave the selection on D4 the way it was when the sheet was
# last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit_Order Number Assign")
$ws.Activate()

# --- bump the seed/default numbers (Last Project Number / Last Purchased
#     Lot Number) from 1 to 10 ---
$ws.Range("A2").Value = 10
$ws.Range("D2").Value = 10

# --- re-fit column widths to the header captions (A, B, C, E); D already
#     carries its own best-fit width from a previous pass and is left as-is ---
$ws.Columns("A:A").ColumnWidth = 16.944010416666668
$ws.Columns("B:B").ColumnWidth = 29.608072916666668
$ws.Columns("C:C").ColumnWidth = 26.053385416666668
$ws.Columns("E:E").ColumnWidth = 25.166666666666668

# --- leave the cursor/selection where the author left it ---
$ws.Range("D4").Select()
